# Update FuelPrices at 2025-03-28 02:23
# - Append a new data row (row 13) to Sheet1 with the latest fuel price entry.
# - The previous "latest" row (row 12) no longer gets the special
#   date-only formatting; it reverts to the standard date/time format,
#   and the newly appended row takes over the date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (previously the last row) switches from the "date only" format
# back to the regular date/time format used by the rest of the column.
$ws.Range("B12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 13 with the latest fuel price record.
$ws.Range("A13").Value = 806.651
$ws.Range("B13").Value = 45734
$ws.Range("B13").NumberFormat = "YYYY-MM-DD"
$ws.Range("C13").Value = 790.4
